# 2026 lab meeting - part 1
# Adds a new "2026" worksheet (after the existing "2025" sheet) and
# populates it with the weekly lab-meeting schedule for Jan-Apr 2026.

$wb = $excel.ActiveWorkbook
$ws2025 = $wb.Worksheets.Item(1)

# --- Create the new "2026" sheet right after "2025" ---------------------
$ws2026 = $wb.Worksheets.Add($null, $ws2025)
$ws2026.Name = "2026"

# --- Pull header + number formatting over from the 2025 sheet -----------
# Header row: Date | Lab Meeting | Journal Club | Food (no nuts) | Cancelled | Notes
# (the 2026 sheet drops the "Technique of the Week" column that 2025 has)
$ws2025.Range("A1:C1").Copy($ws2026.Range("A1:C1"))
$ws2025.Range("E1:G1").Copy($ws2026.Range("D1:F1"))

# Date column format (column A) - reuse the same date-number-format cells
$ws2025.Range("A2:A29").Copy($ws2026.Range("A2:A29"))

# --- Column widths (approximate the target layout) ----------------------
$ws2026.Columns.Item(1).ColumnWidth = 14.02
$ws2026.Columns.Item(2).ColumnWidth = 26.17
$ws2026.Columns.Item(3).ColumnWidth = 18.59
$ws2026.Columns.Item(4).ColumnWidth = 18.02
$ws2026.Columns.Item(5).ColumnWidth = 12.88

# --- Weekly data rows (Jan 12 2026 - Apr 27 2026) ------------------------
# Columns: A=Date B=Lab Meeting C=Journal Club D=Food(no nuts) E=Cancelled F=Notes
$data = @(
  @{ Row=2;  Date=46034; B="lab planning";                 C=$null;                   F="Lab Planning " },
  @{ Row=3;  Date=46041; B="no lab meeting";                C=$null;                   F="MLK Day Volunteering " },
  @{ Row=4;  Date=46048; B="Justin Ma";                     C=$null;                   F="lab farewell party for Justin" },
  @{ Row=5;  Date=46055; B="Guofu Shen";                    C="Daniel Brock (pilot)";  F=$null },
  @{ Row=6;  Date=46062; B="Daniel Brock";                  C="Solomon Gibson";        F="Lab vision / Frankfort lab retreat" },
  @{ Row=7;  Date=46069; B="Solomon Gibson";                C="Chantel George";        F=$null },
  @{ Row=8;  Date=46076; B="Chantel George";                C="Salim Khondker";        F=$null },
  @{ Row=9;  Date=46083; B="Salim Khondker";                C="Soumi Mitra";           F=$null },
  @{ Row=10; Date=46090; B="Soumi Mitra";                   C="Guofu Shen";            F=$null },
  @{ Row=11; Date=46097; B="Ben Frankfort & Daniel Brock";  C=$null;                   F="spring break for Rice" },
  @{ Row=12; Date=46104; B="Med Student Day";                C=$null;                   F="Reagan, Ritu, Miles, Molly, Grace" },
  @{ Row=13; Date=46111; B="Joy Kim";                       C=$null;                   F=$null },
  @{ Row=14; Date=46118; B="no lab meeting";                C=$null;                   F="Ben out of town" },
  @{ Row=15; Date=46125; B="Kevin Wu";                      C=$null;                   F=$null },
  @{ Row=16; Date=46132; B="Giselle Gonzalez";              C=$null;                   F=$null },
  @{ Row=17; Date=46139; B=$null;                           C=$null;                   F="Rice finals" }
)

foreach ($entry in $data) {
  $r = $entry.Row
  $ws2026.Cells.Item($r, 1).Value = $entry.Date
  if ($entry.B) { $ws2026.Cells.Item($r, 2).Value = $entry.B }
  if ($entry.C) { $ws2026.Cells.Item($r, 3).Value = $entry.C }
  $ws2026.Cells.Item($r, 5).Value = $false
  if ($entry.F) { $ws2026.Cells.Item($r, 6).Value = $entry.F }
}

# Remaining rows (18-29) keep the date-column styling but stay blank,
# matching the look of the 2025 sheet's trailing rows.

# --- View settings for the new sheet -------------------------------------
$ws2026.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws2026.Range("D16").Select()

# Make "2026" the active/visible tab (matches tabSelected + activeTab).
$ws2026.Select()
